$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 188; this shifts existing rows 188:291 down to 189:292
$ws.Rows.Item(188).Insert()

# Populate the newly inserted row 188 with the new data record
$ws.Cells.Item(188, 1).Value = 10
$ws.Cells.Item(188, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(188, 3).Value = "La Araucanía"
$ws.Cells.Item(188, 4).Value = 44488
$ws.Cells.Item(188, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(188, 5).Value = 9
$ws.Cells.Item(188, 6).Value = "Fruta"
$ws.Cells.Item(188, 7).Value = 100108
$ws.Cells.Item(188, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(188, 9).Value = 100108005
$ws.Cells.Item(188, 10).Value = "Piña"
$ws.Cells.Item(188, 11).Value = "Caramelo"
$ws.Cells.Item(188, 12).Value = "Segunda"
$ws.Cells.Item(188, 13).Value = 50
$ws.Cells.Item(188, 14).Value = 23000
$ws.Cells.Item(188, 15).Value = 23000
$ws.Cells.Item(188, 16).Value = 23000
$ws.Cells.Item(188, 17).Value = "`$/caja 12 unidades"
$ws.Cells.Item(188, 18).Value = "Ecuador"
$ws.Cells.Item(188, 19).Value = 1917
$ws.Cells.Item(188, 20).Value = 12
